$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three product rows 11-13 get their content rotated up by one (with wraparound):
#   old row 12 content -> row 11
#   old row 13 content -> row 12
#   old row 11 content -> row 13
# (E11 stays blank both before and after, so it is left untouched.)

$ws.Range("A11").Value = "Carta igienica Mini Jumbo OEM fornitore 2400g pacchetto 5 strati di pasta vergine rotoli 165*110mm per ristoranti supermercati"
$ws.Range("B11").Value = "0,2079-0,2426 €"
$ws.Range("C11").Value = "Ordine minimo: 2.000 rulli"
$ws.Range("D11").Value = "Qingdao Shuncai Trading Co., Ltd."

$ws.Range("A12").Value = "Rotoli di carta igienica Mini Jumbo all'ingrosso rotolo di carta velina vergine e Nupkin fabbricazione Maxi miglior toumet Premium confortevole"
$ws.Range("B12").Value = "0,6063-0,6756 €"
$ws.Range("C12").Value = "Ordine minimo: 15.000 rulli"
$ws.Range("D12").Value = "Qingdao Dongfang Jiarui Int'l Co., Ltd."
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.7"

$ws.Range("A13").Value = "Rotolo salute di alta qualità all'ingrosso prezzo di fabbrica Mini rotolo Jumbo carta igienica carta igienica 2 strati miscela pasta di legno 15-19gsm 22-24cm"
$ws.Range("B13").Value = "0,6842-0,8575 €"
$ws.Range("C13").Value = "Ordine minimo: 15.000 rulli"
$ws.Range("D13").Value = "Dongguan Winall Paper Co., Ltd."
$ws.Range("E13").ClearContents()
